# Apply a cyclic rotation of data among rows 20-23 (1-based Excel rows)
# for columns: A, B, D, E, F, G, H, Q, R, Z, AB
#
# Resulting mapping (new row <- old row):
#   row 20 <- row 23
#   row 21 <- row 22
#   row 22 <- row 20
#   row 23 <- row 21

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cols = @("A", "B", "D", "E", "F", "G", "H", "Q", "R", "Z", "AB")

# Capture the current (before) values for each relevant cell in rows 20-23
$values = @{}
foreach ($row in 20..23) {
    foreach ($col in $cols) {
        $values["$col$row"] = $ws.Range("$col$row").Value2
    }
}

# Mapping of new row -> source (old) row
$rowMap = @{
    20 = 23
    21 = 22
    22 = 20
    23 = 21
}

foreach ($newRow in $rowMap.Keys) {
    $oldRow = $rowMap[$newRow]
    foreach ($col in $cols) {
        $ws.Range("$col$newRow").Value2 = $values["$col$oldRow"]
    }
}
